# feat: add 2022-Q1 data
# - Insert a new "2022-Q1" sheet (fund holdings) between "2021-Q4" and "总计"
# - Update the "总计" (totals) summary sheet with a new row for 2022-Q1

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$oldTotal = $wb.Worksheets.Item("总计")

# Capture the old summary rows (date / count / value) before we touch anything.
$oldRows = @(
    @("2021-Q4", 2, 0.59),
    @("2021-Q3", 7, 1.65),
    @("2021-Q2", 4, 0.7),
    @("2021-Q1", 1, 0.01)
)

# Remove the old "总计" sheet; we'll rebuild it after the new quarter sheet
# so the tab order + sheetId allocation come out as: ... 2021-Q4, 2022-Q1, 总计
[void]$oldTotal.Delete()

# --- New "2022-Q1" worksheet -------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Match header/index formatting from the 2021-Q4 sheet.
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'519908"
$q1.Range("C2").Value = "华夏兴华混合A"
$q1.Range("D2").Value = "'9.39"
$q1.Range("E2").Value = "'91.83"
$q1.Range("F2").Value = "'5.87"
$q1.Range("G2").Value = "'0.5512"
$q1.Range("H2").Value = 3

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'960004"
$q1.Range("C3").Value = "华夏兴华混合H"
$q1.Range("D3").Value = "'9.39"
$q1.Range("E3").Value = "'91.83"
$q1.Range("F3").Value = "'5.87"
$q1.Range("G3").Value = "'0.5512"
$q1.Range("H3").Value = 3

# --- Rebuilt "总计" worksheet -------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$q1.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$q1.Range("A2:A3").Copy()
$total.Range("A2:A6").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 1.1

for ($i = 0; $i -lt $oldRows.Length; $i++) {
    $r = $i + 3
    $row = $oldRows[$i]
    $total.Range("A$r").Value = $i + 1
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
}
